$d = $word.ActiveDocument

# The new row goes at the end of the last table ("Type contact"),
# right after the existing "type" row.
$table = $d.Tables.Item($d.Tables.Count)

$newRow = $table.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "details"
$newRow.Cells.Item(2).Range.Text = "Détails de contact"
$newRow.Cells.Item(3).Range.Text = "string"
$newRow.Cells.Item(4).Range.Text = "0..1"
$newRow.Cells.Item(5).Range.Text = "1. RFGI (si RFGI disponible)" + [char]11 + "2. Numéro de téléphone"
$newRow.Cells.Item(6).Range.Text = "0612342536"
